$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.502.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.806.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.582"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.22%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.50%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0954"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.799.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.478.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0800"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.684"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.396.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.954"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.966.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("E51").Value = "  +0.14%  "
